$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout:
#   Row 1 (header): B1:F1 = PercActivationsOriginal/Correct/Incorrect/Corrected/Fixed
#   Rows 2..N:      A = segment name, B:F = the five percentage metrics
#
# Target layout:
#   Row 1 (header): B1 = "segments", C1:G1 = the five original headers (shifted right)
#   Rows 2..N:      A = numeric index (0-based), B = segment name, C:G = the metrics

$lastRow = $ws.UsedRange.Rows.Count

# Shift everything one column to the right (A->B, B->C, ... E->F, F->G),
# leaving a blank column A in its place.
$ws.Columns.Item(1).Insert()

# The (now shifted) segment-name column landed in B with the bordered/bold
# header style carried over from the old column A - the target data rows
# have no explicit style there, so clear it back to the default.
$ws.Range("B2:B" + $lastRow).ClearFormats()

# Give the new header cell (B1) and the new index column (A2:A<lastRow>) the
# same bordered/bold/centered style used by the rest of row 1 (style index
# "1" in the original file) by copying formats from an existing styled cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("A2:A" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header text.
$ws.Range("B1").Value = "segments"

# New zero-based numeric index column.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
